# Apply the changes described by the commit "Added data and V2 File":
#  - Header cell A1 changes from "Trend" to "Year"
#  - Active selection on the sheet moves from G13 to A2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header text in A1.
$ws.Range("A1").Value = "Year"

# Move / record the active cell selection to A2, matching the saved view state.
[void]$ws.Range("A2").Select()
